$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition the window / active selection (cosmetic, matches author's session) ---
$excel.Left = 3780
$excel.Top = 2980

# --- Make room for the new "assocId" rows above each existing association's id row ---
$ws.Rows("47:47").Insert()
$ws.Rows("50:50").Insert()
$ws.Rows("53:53").Insert()
$ws.Rows("56:56").Insert()
$ws.Rows("59:59").Insert()
$ws.Rows("62:62").Insert()
$ws.Rows("65:65").Insert()
$ws.Rows("67:67").Insert()

# --- Make room for three brand-new "enemyOf" / "enemyOF" associations ---
$ws.Rows("71:76").Insert()

# --- The sheet only grew by 10 rows net (reuse some of the old blank filler rows) ---
$ws.Rows("84:87").Delete()

# --- Fill in the new assocId rows that precede each existing association block ---
$ws.Range("A47").Value = "assocId"
$ws.Range("B47").Value = "23"

$ws.Range("A50").Value = "assocId"
$ws.Range("B50").Value = "65"

$ws.Range("A53").Value = "assocId"
$ws.Range("B53").Value = "54"

$ws.Range("A56").Value = "assocId"
$ws.Range("B56").Value = "96"

$ws.Range("A59").Value = "assocId"
$ws.Range("B59").Value = "21"

$ws.Range("A62").Value = "assocId"
$ws.Range("B62").Value = "38"

$ws.Range("A65").Value = "assocId"
$ws.Range("B65").Value = "42"

$ws.Range("A67").Value = "assocId"
$ws.Range("B67").Value = "55"

# --- Fill in the three brand-new enemyOf / enemyOF association blocks ---
# (filled in the same order the original author typed them: the two
# "enemyOf"/"enemyOF" template rows first, then the remaining assocId rows)
$ws.Range("A71").Value = "assocId"
$ws.Range("B71").Value = "74"

$ws.Range("A72").Value = "220"
$ws.Range("B72").Value = "enemyOf"
$ws.Range("C72").Value = "001"

$ws.Range("A76").Value = "A15"
$ws.Range("B76").Value = "enemyOF"
$ws.Range("C76").Value = "002"

$ws.Range("A73").Value = "assocId"
$ws.Range("B73").Value = "12"

$ws.Range("A74").Value = "220"
$ws.Range("B74").Value = "enemyOf"
$ws.Range("C74").Value = "777"

$ws.Range("A75").Value = "assocId"
$ws.Range("B75").Value = "15"

# --- Restore the author's viewport / selection on the sheet ---
$ws.Application.Goto($ws.Range("A23"), $false)
$ws.Range("D72").Select()
